$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-8 with re-recorded values
$ws.Range("C2").Value = 5.085351
$ws.Range("D2").Value = 0.000000

$ws.Range("B3").Value = -0.007835
$ws.Range("C3").Value = 5.085351

$ws.Range("B4").Value = -0.012120
$ws.Range("C4").Value = 5.085351

$ws.Range("B5").Value = -0.016405
$ws.Range("C5").Value = 5.085351

$ws.Range("B6").Value = -0.020690
$ws.Range("C6").Value = 5.085351

$ws.Range("B7").Value = -0.024975
$ws.Range("C7").Value = 5.232569
$ws.Range("D7").Value = 0.147217

$ws.Range("B8").Value = -0.029260
$ws.Range("C8").Value = 5.821439
$ws.Range("D8").Value = 0.736087

# Remove row 9 entirely (was: 7.000000, 0.026445, -6.201318, 1.079594)
$ws.Range("A9:D9").Delete()
